$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 871.625
$ws.Cells.Item(92, 9).Value = 871.625
$ws.Cells.Item(92, 11).Value = 871.625
$ws.Cells.Item(92, 13).Value = 376.375
$ws.Cells.Item(94, 8).Value = 5564.091
$ws.Cells.Item(96, 8).Value = 1527.3334
$ws.Cells.Item(96, 9).Value = 1747.6666
$ws.Cells.Item(96, 10).Value = 866.3333
$ws.Cells.Item(96, 11).Value = 5242.9998
$ws.Cells.Item(96, 12).Value = 2598.9999
$ws.Cells.Item(96, 13).Value = -3869.9998
$ws.Cells.Item(96, 14).Value = -5344.9999
$ws.Cells.Item(97, 8).Value = 5499.5
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 5499.5
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 16498.5
$ws.Cells.Item(97, 13).ClearContents()
$ws.Cells.Item(97, 14).Value = -17490.5
$ws.Cells.Item(100, 8).Value = 3026.8948
$ws.Cells.Item(100, 9).Value = 2654.2307
$ws.Cells.Item(100, 11).Value = 2654.2307
$ws.Cells.Item(100, 13).Value = -2113.2307
$ws.Cells.Item(101, 8).Value = 974.7273
$ws.Cells.Item(101, 9).Value = 554.75
$ws.Cells.Item(101, 10).Value = 2094.6667
$ws.Cells.Item(101, 11).Value = 1664.25
$ws.Cells.Item(101, 12).Value = 6284.000100000001
$ws.Cells.Item(101, 13).Value = -42.25
$ws.Cells.Item(101, 14).Value = -9528.000100000001
$ws.Cells.Item(104, 8).Value = 156.4
$ws.Cells.Item(104, 9).Value = 156.4
$ws.Cells.Item(104, 11).Value = 469.2
$ws.Cells.Item(104, 13).Value = 1277.8
$ws.Cells.Item(112, 8).Value = 2678.6843
$ws.Cells.Item(112, 10).Value = 2799.7222
$ws.Cells.Item(112, 12).Value = 8399.1666
$ws.Cells.Item(112, 14).Value = -10615.1666
$ws.Cells.Item(115, 8).Value = 1186.1111
$ws.Cells.Item(115, 9).Value = 711.875
$ws.Cells.Item(115, 10).Value = 4980
$ws.Cells.Item(115, 11).Value = 2135.625
$ws.Cells.Item(115, 12).Value = 14940
$ws.Cells.Item(115, 13).Value = -568.625
$ws.Cells.Item(115, 14).Value = -18074
$ws.Cells.Item(129, 8).Value = 1074.127
$ws.Cells.Item(129, 9).Value = 424
$ws.Cells.Item(129, 10).Value = 1211.6538
$ws.Cells.Item(129, 11).Value = 1272
$ws.Cells.Item(129, 12).Value = 3634.9614
$ws.Cells.Item(129, 13).Value = 3728
$ws.Cells.Item(129, 14).Value = -13634.9614
$ws.Cells.Item(137, 8).Value = 1521.9546
$ws.Cells.Item(137, 9).Value = 1392
$ws.Cells.Item(137, 10).Value = 1800.4286
$ws.Cells.Item(137, 11).Value = 4176
$ws.Cells.Item(137, 12).Value = 5401.2858
$ws.Cells.Item(137, 13).Value = -1626
$ws.Cells.Item(137, 14).Value = -10501.2858
$ws.Cells.Item(138, 8).Value = 2508.6667
$ws.Cells.Item(138, 9).Value = 1450.175
$ws.Cells.Item(138, 10).Value = 4349.522
$ws.Cells.Item(138, 11).Value = 4350.525
$ws.Cells.Item(138, 12).Value = 13048.566
$ws.Cells.Item(138, 13).Value = 789.4750000000004
$ws.Cells.Item(138, 14).Value = -23328.566

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11883.566
$ws.Cells.Item(32, 9).Value = 11776.7705
$ws.Cells.Item(32, 10).Value = 12908.8
$ws.Cells.Item(32, 11).Value = 11776.7705
$ws.Cells.Item(32, 12).Value = 12908.8
$ws.Cells.Item(32, 13).Value = -11489.7705
$ws.Cells.Item(32, 14).Value = -13482.8
$ws.Cells.Item(97, 8).Value = 949.8889
$ws.Cells.Item(97, 9).Value = 633.2143
$ws.Cells.Item(97, 10).Value = 2058.25
$ws.Cells.Item(97, 11).Value = 633.2143
$ws.Cells.Item(97, 12).Value = 2058.25
$ws.Cells.Item(97, 13).Value = -137.2143
$ws.Cells.Item(97, 14).Value = -3050.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2486.926
$ws.Cells.Item(31, 9).Value = 1574.8334
$ws.Cells.Item(31, 10).Value = 4311.1113
$ws.Cells.Item(31, 11).Value = 1574.8334
$ws.Cells.Item(31, 12).Value = 4311.1113
$ws.Cells.Item(31, 13).Value = -1279.8334
$ws.Cells.Item(31, 14).Value = -4901.1113
$ws.Cells.Item(34, 8).Value = 2486.926
$ws.Cells.Item(34, 9).Value = 1574.8334
$ws.Cells.Item(34, 10).Value = 4311.1113
$ws.Cells.Item(34, 11).Value = 1574.8334
$ws.Cells.Item(34, 12).Value = 4311.1113
$ws.Cells.Item(34, 13).Value = -1372.8334
$ws.Cells.Item(34, 14).Value = -4715.1113
$ws.Cells.Item(105, 8).Value = 966.8333
$ws.Cells.Item(105, 9).Value = 900.5
$ws.Cells.Item(105, 10).Value = 1000
$ws.Cells.Item(105, 11).Value = 900.5
$ws.Cells.Item(105, 12).Value = 1000
$ws.Cells.Item(105, 13).Value = 846.5
$ws.Cells.Item(105, 14).Value = -4494
$ws.Cells.Item(141, 8).Value = 39536.625
$ws.Cells.Item(141, 10).Value = 39536.625
$ws.Cells.Item(141, 12).Value = 39536.625
$ws.Cells.Item(141, 14).Value = -49896.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2069.0715
$ws.Cells.Item(5, 9).Value = 2807.75
$ws.Cells.Item(5, 10).Value = 1084.1666
$ws.Cells.Item(5, 11).Value = 8423.25
$ws.Cells.Item(5, 12).Value = 3252.4998
$ws.Cells.Item(5, 13).Value = -8311.25
$ws.Cells.Item(5, 14).Value = -3476.4998
$ws.Cells.Item(122, 8).Value = 881.1429000000001
$ws.Cells.Item(122, 9).Value = 660
$ws.Cells.Item(122, 10).Value = 1047
$ws.Cells.Item(122, 11).Value = 5940
$ws.Cells.Item(122, 12).Value = 9423
$ws.Cells.Item(122, 13).Value = -3490
$ws.Cells.Item(122, 14).Value = -14323
$ws.Cells.Item(135, 8).Value = 2069.0715
$ws.Cells.Item(135, 9).Value = 2807.75
$ws.Cells.Item(135, 10).Value = 1084.1666
$ws.Cells.Item(135, 11).Value = 25269.75
$ws.Cells.Item(135, 12).Value = 9757.499400000001
$ws.Cells.Item(135, 13).Value = -22734.75
$ws.Cells.Item(135, 14).Value = -14827.4994

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 150680.1
$ws.Cells.Item(97, 9).Value = 84001.664
$ws.Cells.Item(97, 10).Value = 250697.75
$ws.Cells.Item(97, 11).Value = 84001.664
$ws.Cells.Item(97, 12).Value = 250697.75
$ws.Cells.Item(97, 13).Value = -83505.664
$ws.Cells.Item(97, 14).Value = -251689.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1971.8572
$ws.Cells.Item(93, 9).Value = 1398.6666
$ws.Cells.Item(93, 10).Value = 2401.75
$ws.Cells.Item(93, 11).Value = 1398.6666
$ws.Cells.Item(93, 12).Value = 2401.75
$ws.Cells.Item(93, 13).Value = -150.6666
$ws.Cells.Item(93, 14).Value = -4897.75
$ws.Cells.Item(100, 8).Value = 8737.143
$ws.Cells.Item(100, 9).Value = 17320
$ws.Cells.Item(100, 10).Value = 2300
$ws.Cells.Item(100, 11).Value = 17320
$ws.Cells.Item(100, 12).Value = 2300
$ws.Cells.Item(100, 13).Value = -16779
$ws.Cells.Item(100, 14).Value = -3382

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1874.5
$ws.Cells.Item(96, 9).Value = 1311.5
$ws.Cells.Item(96, 11).Value = 1311.5
$ws.Cells.Item(96, 13).Value = 61.5
$ws.Cells.Item(100, 8).Value = 7700.2856
$ws.Cells.Item(100, 9).Value = 11549.889
$ws.Cells.Item(100, 10).Value = 771
$ws.Cells.Item(100, 11).Value = 23099.778
$ws.Cells.Item(100, 12).Value = 1542
$ws.Cells.Item(100, 13).Value = -22558.778
$ws.Cells.Item(100, 14).Value = -2624
$ws.Cells.Item(136, 8).Value = 1401.8529
$ws.Cells.Item(136, 9).Value = 1482.56
$ws.Cells.Item(136, 10).Value = 1177.6666
$ws.Cells.Item(136, 11).Value = 4447.68
$ws.Cells.Item(136, 12).Value = 3532.9998
$ws.Cells.Item(136, 13).Value = -1897.68
$ws.Cells.Item(136, 14).Value = -8632.9998
